$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Format")

# Row 4: "Big font" label - merge A4:C4, centered
$ws.Range("A4:C4").HorizontalAlignment = -4108
$ws.Range("A4:C4").Merge()

# Row 5: "Medium font" label - merge A5:C5, centered
$ws.Range("A5:C5").HorizontalAlignment = -4108
$ws.Range("A5:C5").Merge()

# Row 6: "Big, medium, and small fonts." - merge A6:E6, centered
$ws.Range("A6:E6").HorizontalAlignment = -4108
$ws.Range("A6:E6").Merge()

# Row 13: color name labels for background colors demo
$ws.Range("A13").Value = "red"
$ws.Range("B13").Value = "yellow"
$ws.Range("C13").Value = "blue"
$ws.Range("D13").Value = "purple"
$ws.Range("E13").Value = "light green"

# Activate this sheet and move the selection, matching the saved view state
$ws.Activate()
[void]$ws.Range("G4").Select()
